$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D, M, N, O, P, R, S
# (result of shuffling the data among rows 2,4,6,7,8,9,10,11)
$data = @{
    2  = @{ D = 44567; M = 80;  N = 2400; O = 2400; P = 2400; R = "Región de La Araucanía"; S = 2400 }
    4  = @{ D = 44176; M = 20;  N = 3000; O = 3000; P = 3000; R = "Región de O'Higgins";    S = 3000 }
    6  = @{ D = 44323; M = 20;  N = 3200; O = 3200; P = 3200; R = "Región de La Araucanía"; S = 3200 }
    7  = @{ D = 44214; M = 50;  N = 1800; O = 1800; P = 1800; R = "Región de La Araucanía"; S = 1800 }
    8  = @{ D = 44215; M = 65;  N = 2800; O = 2800; P = 2800; R = "Región de La Araucanía"; S = 2800 }
    9  = @{ D = 44592; M = 5;   N = 7500; O = 7500; P = 7500; R = "Región de La Araucanía"; S = 7500 }
    10 = @{ D = 44551; M = 120; N = 4500; O = 4500; P = 4500; R = "Región de O'Higgins";    S = 4500 }
    11 = @{ D = 44616; M = 200; N = 3200; O = 3200; P = 3200; R = "Región de La Araucanía"; S = 3200 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
